$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark three more "corregido" (fixed) entries in column B
$ws.Range("B10").Value = "corregido"
$ws.Range("B11").Value = "corregido"
$ws.Range("B12").Value = "corregido"

# Fill the previously-empty row 15 with the text that used to live in A16
$ws.Range("A15").Value = "al apretar Partidos en Vivo, entra en todos los partidos"

# New feedback items appended after the existing list
$ws.Range("A16").Value = "falta un salir de la visualizacion del partido en espectador"
$ws.Range("A17").Value = 'falta un "suspender" partido para que el partido quede guardado todos los datos al momento de la suspension '
$ws.Range("A18").Value = "al momento de agregar una falta, deberia abrirse un menu flotante que indique que tipo de falta es (personal, tecnica, antideportiva, descalificadora)"
$ws.Range("A19").Value = "en caso de sumar 2 faltas tecnicas, 2 faltas antideportivas o 1 falta tecn y 1 desc, el jugador queda descalificado de ese partido debe aparecer un GD"
$ws.Range("A20").Value = "el jugador que hace 5 faltas esta bien que no pueda tener mas acciones, pero debe poder hacer sustitucion para que entre otro jugador suplente"
$ws.Range("A22").Value = "al momento de cargar jugadores nuevos se tiene que poder marcar al mismo como jugador Refuerzo, el cual va a tener reglas especiales en cuanto a la cantidad de cuartos que puede jugar"
$ws.Range("A21").Value = "al iniciar el juego, pero antes de elegir a los 5 titulares, los dos equipos deberian poder elegir a los 12 jugadores citados para el juego"

# Move the active selection to follow the newly-entered data
$ws.Range("A23").Select()
